# feat: quick look in progress. complete loading modal window in object
#
# Adds the "quick look" modal locators to the PDP_Page sheet:
#   - a modal-scoped "original price" xpath (new row under the existing one)
#   - a modal-scoped "sale price" xpath (new row under the existing one)
#   - fixes the "quick look popul" typo -> "quick look popup"
#   - adds a new "sku" locator row
#   - a couple of trailing blank rows to match the extended table

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("PDP_Page")

# --- insert the modal "original price" row right after the existing one (row 7) ---
$ws.Rows.Item(8).Insert()
$ws.Range("A8").Value = "original price"
$ws.Range("B8").Value = "xpath"
$ws.Range("C8").Value = '//*[@id="modal"]/div[1]/div/section[2]/div[1]/div/div[1]/span'

# --- insert the modal "sale price" row right after the existing sale price row (now row 9) ---
$ws.Rows.Item(10).Insert()
$ws.Range("A10").Value = "sale price"
$ws.Range("B10").Value = "xpath"
$ws.Range("C10").Value = '//*[@id="modal"]/div[1]/div/section[2]/div[1]/div/div[2]/span[2]'

# --- fill in the next (already-blank, pre-formatted) row with the new "sku" locator ---
$ws.Range("B19").Value = "xpath"
$ws.Range("C19").Value = '//*[@id="content1"]/div[1]/p[2]'

# --- fix the "quick look popul" typo (now shifted down to row 18) ---
$ws.Range("A18").Value = "quick look popup"

$ws.Range("A19").Value = "sku"

# --- update the sheet selection to match ---
$ws.Range("A19").Select() | Out-Null
